$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text so numeric-looking strings (e.g. "43.049.45", "0.638")
# are not auto-converted to numbers by Excel, matching the source data (always text).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.049.45"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.391.85"
$ws.Range("E3").Value = "  +4.79%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "334.67"
$ws.Range("E5").Value = "  +8.34%  "
$ws.Range("D6").Value = "102.75"
$ws.Range("E6").Value = "  -8.56%  "
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.638"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("D10").Value = "41.42"
$ws.Range("E10").Value = "  -6.36%  "
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "8.61"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "1.04"
$ws.Range("E13").Value = "  -4.58%  "
$ws.Range("D14").Value = "16.92"
$ws.Range("E14").Value = "  +8.81%  "
$ws.Range("D15").Value = "0.106"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "2.749.26"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "2.394.54"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "43.015.42"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "3.85"
$ws.Range("E21").Value = "  +6.41%  "
$ws.Range("D22").Value = "77.07"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "272.79"
$ws.Range("E23").Value = "  +6.63%  "
$ws.Range("E24").Value = "  -3.32%  "
$ws.Range("D25").Value = "9.80"
$ws.Range("E25").Value = "  +9.50%  "
$ws.Range("D26").Value = "11.83"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "24.36"
$ws.Range("E28").Value = "  +10.01%  "
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").Value = "174.29"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "3.14"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "36.55"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("D33").Value = "0.0923"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "6.11"
$ws.Range("E34").Value = "  +7.12%  "
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("D36").Value = "4.80"
$ws.Range("E36").Value = "  -4.06%  "
$ws.Range("D37").Value = "4.03"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("D38").Value = "0.0364"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("D40").Value = "2.86"
$ws.Range("E40").Value = "  +12.32%  "
$ws.Range("D41").Value = "1.54"
$ws.Range("E41").Value = "  +12.18%  "
$ws.Range("D42").Value = "0.234"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").Value = "70.00"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "91.51"
$ws.Range("E45").Value = "  +45.94%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "117.43"
$ws.Range("E46").Value = "  +8.51%  "
$ws.Range("D47").Value = "12.13"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "5.54"
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.10"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "0.502"
$ws.Range("E50").Value = "  +14.71%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.613.13"
$ws.Range("E51").Value = "  +8.92%  "

# Restore default style (clears the temporary text-format style so cells
# keep their original unstyled appearance, only the stored type changes to text).
$ws.Range("D2:E51").Style = "Normal"
